$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Participants")

# Fill in participant 8's (row 9) previously-blank record.
# Copy cell formatting from the row above (participant 7, row 8) so the
# green "piloted" styling + number formats match the rest of rows 2-8.
$srcRow = $ws.Range("A8:F8")
$dstRow = $ws.Range("A9:F9")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)  # xlPasteFormats

$srcG = $ws.Range("G8")
$dstG = $ws.Range("G9")
$srcG.Copy()
$dstG.PasteSpecial(-4122)

$srcH = $ws.Range("H8")
$dstH = $ws.Range("H9")
$srcH.Copy()
$dstH.PasteSpecial(-4122)

$ws.Cells.Item(9,2).Value = "Lucia Gomez Lopez"
$ws.Cells.Item(9,3).Value = 11000
$ws.Cells.Item(9,7).Value = "C01"
$ws.Cells.Item(9,8).Value = 45251
$ws.Cells.Item(9,9).Value = "-"

# Fix participant 2's code: "20201" -> "020201" (missing leading zero)
$ws.Cells.Item(3,4).Value = "020201"

# Rows 10-49: fill in the Language_test (column G) values, cycling
# through the four testing locations in step with the existing List
# (column F) groupings of four rows each.
$pattern = @("C02","C03","C04","C01")
for ($r = 10; $r -le 49; $r++) {
  $idx = ($r - 10) % 4
  $ws.Cells.Item($r,7).Value = $pattern[$idx]
}

# Move the active selection on the Participants sheet to I12.
[void]$ws.Select()
[void]$ws.Range("I12").Select()
